# Auto-generated edit script applying the cryptos.xlsx diff
# (crypto price/volume refresh + two row-pair swaps: rows 33/34 and rows 50/51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.364.68"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "3.595.27"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("E4").Value = "  +0.01%  "
$c = $ws.Range("D5")
$origStyle = $c.Style
$c.Value = "'606.89"
$c.Style = $origStyle
$ws.Range("E5").Value = "  +0.17%  "
$c = $ws.Range("D6")
$origStyle = $c.Style
$c.Value = "'148.48"
$c.Style = $origStyle
$ws.Range("E6").Value = "  +2.81%  "
$ws.Range("D7").Value = "3.593.25"
$ws.Range("E7").Value = "  +0.89%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  -0.30%  "
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("E11").Value = "  +0.74%  "
$c = $ws.Range("D12")
$origStyle = $c.Style
$c.Value = "'0.414"
$c.Style = $origStyle
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("D13").Value = "4.202.51"
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("E14").Value = "  -0.43%  "
$c = $ws.Range("D15")
$origStyle = $c.Style
$c.Value = "'29.56"
$c.Style = $origStyle
$ws.Range("E15").Value = "  -1.71%  "
$ws.Range("D16").Value = "3.587.39"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("E17").Value = "  +2.12%  "
$ws.Range("D18").Value = "66.438.13"
$ws.Range("E18").Value = "  +0.33%  "
$c = $ws.Range("D19")
$origStyle = $c.Style
$c.Value = "'11.13"
$c.Style = $origStyle
$ws.Range("E19").Value = "  -2.85%  "
$ws.Range("E20").Value = "  +1.97%  "
$ws.Range("E21").Value = "  +1.41%  "
$c = $ws.Range("D22")
$origStyle = $c.Style
$c.Value = "'423.25"
$c.Style = $origStyle
$ws.Range("E22").Value = "  -1.83%  "
$ws.Range("E23").Value = "  +0.39%  "
$c = $ws.Range("D24")
$origStyle = $c.Style
$c.Value = "'78.72"
$c.Style = $origStyle
$ws.Range("E24").Value = "  -1.34%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").Value = "  +2.63%  "
$c = $ws.Range("D27")
$origStyle = $c.Style
$c.Value = "'8.25"
$c.Style = $origStyle
$ws.Range("E27").Value = "  +4.79%  "
$c = $ws.Range("D28")
$origStyle = $c.Style
$c.Value = "'9.38"
$c.Style = $origStyle
$ws.Range("E28").Value = "  +2.88%  "
$ws.Range("E29").Value = "  -0.19%  "
$c = $ws.Range("D30")
$origStyle = $c.Style
$c.Value = "'0.999"
$c.Style = $origStyle
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").Value = "3.591.29"
$ws.Range("E31").Value = "  +0.80%  "
$ws.Range("E32").Value = "  +3.45%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Range("D33")
$origStyle = $c.Style
$c.Value = "'1.44"
$c.Style = $origStyle
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D34")
$origStyle = $c.Style
$c.Value = "'25.12"
$c.Style = $origStyle
$ws.Range("E34").Value = "  -1.39%  "
$ws.Range("E36").Value = "  -0.55%  "
$c = $ws.Range("D37")
$origStyle = $c.Style
$c.Value = "'5.58"
$c.Style = $origStyle
$ws.Range("E37").Value = "  +0.66%  "
$ws.Range("E38").Value = "  -2.62%  "
$c = $ws.Range("D39")
$origStyle = $c.Style
$c.Value = "'174.91"
$c.Style = $origStyle
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("E41").Value = "  +0.20%  "
$c = $ws.Range("D42")
$origStyle = $c.Style
$c.Value = "'0.885"
$c.Style = $origStyle
$ws.Range("E42").Value = "  -1.01%  "
$c = $ws.Range("D43")
$origStyle = $c.Style
$c.Value = "'45.91"
$c.Style = $origStyle
$ws.Range("E43").Value = "  -0.12%  "
$c = $ws.Range("D44")
$origStyle = $c.Style
$c.Value = "'1.86"
$c.Style = $origStyle
$ws.Range("E44").Value = "  -4.51%  "
$c = $ws.Range("D45")
$origStyle = $c.Style
$c.Value = "'1.00"
$c.Style = $origStyle
$ws.Range("E45").Value = "  +0.02%  "
$c = $ws.Range("D46")
$origStyle = $c.Style
$c.Value = "'2.53"
$c.Style = $origStyle
$ws.Range("E46").Value = "  +4.88%  "
$c = $ws.Range("D47")
$origStyle = $c.Style
$c.Value = "'23.71"
$c.Style = $origStyle
$ws.Range("E47").Value = "  +2.95%  "
$ws.Range("E48").Value = "  +0.34%  "
$c = $ws.Range("D49")
$origStyle = $c.Style
$c.Value = "'24.20"
$c.Style = $origStyle
$ws.Range("E49").Value = "  -3.23%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$c = $ws.Range("D50")
$origStyle = $c.Style
$c.Value = "'1.13"
$c.Style = $origStyle
$ws.Range("E50").Value = "  -5.41%  "
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$c = $ws.Range("D51")
$origStyle = $c.Style
$c.Value = "'0.960"
$c.Style = $origStyle
$ws.Range("E51").Value = "  +3.24%  "
